# Auto-generated Excel COM-interop script to apply the weekly CompStat data refresh
# (new crime data collected for week of 4/15/2024 - 4/21/2024, Volume 31 Number 16)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/issue number and the reporting week dates ---
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# --- Cells that change style category: was shown as text "0"/"***.*" (no data), now a real number ---
# Donor cells already carrying the target numeric style are used to clone exact formatting.
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("H15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("F15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("F15").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = 1
$ws.Range("H15").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = -100
$ws.Range("F15").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("H15").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$ws.Range("H33").Value = -100

# --- Cells that change style category: was a real number, now shown as text "0"/"***.*" (no data) ---
# Toggle to text format so the numeric-looking string is kept literal, write it, then restore the
# donors General-format styling (cloned from an existing text-styled cell) so the final style matches.
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("F31").PasteSpecial(-4122)

# --- Plain value updates (style unchanged) for the weekly crime-stat table, rows 14-33 ---
# Row 14
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -95.833333333333
# Row 15
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 12
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = -60
$ws.Range("M15").Value = -20
$ws.Range("N15").Value = -65.217391304347
# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -46.153846153846
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 33
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 144
$ws.Range("J16").Value = 164
$ws.Range("K16").Value = -12.195121951219
$ws.Range("L16").Value = -21.311475409836
$ws.Range("M16").Value = 6.666666666666
$ws.Range("N16").Value = -76.585365853658
# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = 15.384615384615
$ws.Range("F17").Value = 46
$ws.Range("G17").Value = 65
$ws.Range("H17").Value = -29.230769230769
$ws.Range("I17").Value = 196
$ws.Range("J17").Value = 241
$ws.Range("K17").Value = -18.672199170124
$ws.Range("L17").Value = -2.487562189054
$ws.Range("M17").Value = 28.104575163398
$ws.Range("N17").Value = -20
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -46.153846153846
$ws.Range("I18").Value = 90
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = -24.369747899159
$ws.Range("L18").Value = 3.448275862068
$ws.Range("M18").Value = -21.739130434782
$ws.Range("N18").Value = -82.558139534883
# Row 19
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -19.047619047619
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = -23.809523809523
$ws.Range("I19").Value = 265
$ws.Range("J19").Value = 295
$ws.Range("K19").Value = -10.169491525423
$ws.Range("L19").Value = -16.403785488959
$ws.Range("M19").Value = 75.496688741721
$ws.Range("N19").Value = 16.740088105726
# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -46.153846153846
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 58
$ws.Range("H20").Value = -51.724137931034
$ws.Range("I20").Value = 136
$ws.Range("J20").Value = 204
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -23.163841807909
$ws.Range("M20").Value = 115.873015873016
$ws.Range("N20").Value = -77.257525083612
# Row 21
$ws.Range("C21").Value = 50
$ws.Range("D21").Value = 71
$ws.Range("E21").Value = -29.577464788732
$ws.Range("F21").Value = 189
$ws.Range("G21").Value = 270
$ws.Range("H21").Value = -30
$ws.Range("I21").Value = 840
$ws.Range("J21").Value = 1038
$ws.Range("K21").Value = -19.075144508670
$ws.Range("L21").Value = -15.065722952477
$ws.Range("M21").Value = 33.545310015898
$ws.Range("N21").Value = -62.633451957295
# Row 22
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = -25
# Row 23
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 76
$ws.Range("J23").Value = 88
$ws.Range("K23").Value = -13.636363636363
$ws.Range("L23").Value = -26.923076923076
$ws.Range("M23").Value = 13.432835820895
# Row 24
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -45.714285714285
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = -27.611940298507
$ws.Range("I24").Value = 517
$ws.Range("J24").Value = 609
$ws.Range("K24").Value = -15.106732348111
$ws.Range("L24").Value = -5.137614678899
$ws.Range("M24").Value = 28.606965174129
# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -70.588235294117
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 55
$ws.Range("H25").Value = -32.727272727272
$ws.Range("I25").Value = 184
$ws.Range("J25").Value = 266
$ws.Range("K25").Value = -30.827067669172
$ws.Range("L25").Value = -43.558282208589
# Row 26
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 73
$ws.Range("G26").Value = 74
$ws.Range("H26").Value = -1.351351351351
$ws.Range("I26").Value = 350
$ws.Range("J26").Value = 314
$ws.Range("K26").Value = 11.464968152866
$ws.Range("L26").Value = 6.060606060606
$ws.Range("M26").Value = -21.171171171171
# Row 27
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = -44.827586206896
# Row 28
$ws.Range("C28").Value = 6
$ws.Range("F28").Value = 11
$ws.Range("H28").Value = 175
$ws.Range("I28").Value = 37
$ws.Range("K28").Value = 12.121212121212
$ws.Range("L28").Value = 94.736842105263
# Row 29
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -30.769230769230
$ws.Range("M29").Value = -25
$ws.Range("N29").Value = -79.069767441860
# Row 30
$ws.Range("H30").Value = -100
$ws.Range("L30").Value = -38.461538461538
$ws.Range("M30").Value = -20
$ws.Range("N30").Value = -80
# Row 33
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = -50

$wb.Application.CutCopyMode = $false

